# "now adds shopping list" - append a trailing comma to the
# fresh/shopping-list ingredient cells so they read as a proper
# comma-separated shopping list (matches the other already-comma-terminated
# ingredient cells in the sheet).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B12").Value = "rice, coconut milk,"
$ws.Range("B17").Value = "burgers, fries,"
$ws.Range("B18").Value = "flour, egg, milk,"
$ws.Range("B4").Value  = "fiskeboller, potato, carrots,"
$ws.Range("D2").Value  = "paprika, carrot, potato, onion,"
$ws.Range("D3").Value  = "carrot, potato,"
$ws.Range("D10").Value = "potato,"
$ws.Range("D12").Value = "carrot,"
$ws.Range("D16").Value = "carrot, spring onion,"

# Leave the cursor where the author ended up when they made this edit.
$ws.Range("D17").Select() | Out-Null
